# Updated symbol list: refresh price (D) and a couple of volume-label (E) cells
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.31"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "20.92"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.219"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06203"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.579"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.554"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.485"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8224"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01386"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1641"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08281"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03512"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03101"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09126"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.777"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001642"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04689"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006436"
$ws.Range("E19").Value = "18TigerCashTCHBestin24h"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006159"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.822"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.359"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3305"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1250"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002738"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04679"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007019"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003523"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01087"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006187"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8456"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001968"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00001901"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01241"

$wb.Save()